$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C to fit updated content
$ws.Columns.Item(3).ColumnWidth = 54.65

# Add a total row summing the Qty column (D2:D7)
$ws.Range("D8").Formula = "=SUM(D2:D7)"

# Update the active selection to reflect where the user ended up (D9)
$ws.Range("D9").Select()
